$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 79.5
$ws.Range("I9").Value = 59.95
$ws.Range("J9").Value = 275
$ws.Range("K9").Value = 59.95
$ws.Range("L9").Value = 275
$ws.Range("M9").Value = 109.05
$ws.Range("N9").Value = -613

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1441.3077
$ws.Range("I38").Value = 108.14286
$ws.Range("J38").Value = 2996.6667
$ws.Range("K38").Value = 324.42858
$ws.Range("L38").Value = 8990.000100000001
$ws.Range("M38").Value = 47.57141999999999
$ws.Range("N38").Value = -9734.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1106.4286
$ws.Range("I58").Value = 829
$ws.Range("K58").Value = 2487
$ws.Range("M58").Value = -2337

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 33977
$ws.Range("J95").Value = 33977
$ws.Range("L95").Value = 33977
$ws.Range("N95").Value = -39469

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2551.0334
$ws.Range("I138").Value = 2176.125
$ws.Range("J138").Value = 2757.8794
$ws.Range("K138").Value = 6528.375
$ws.Range("L138").Value = 8273.638199999999
$ws.Range("M138").Value = -1388.375
$ws.Range("N138").Value = -18553.6382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31341.684
$ws.Range("I32").Value = 32114.07
$ws.Range("J32").Value = 16666.334
$ws.Range("K32").Value = 32114.07
$ws.Range("L32").Value = 16666.334
$ws.Range("M32").Value = -31827.07
$ws.Range("N32").Value = -17240.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2112.75
$ws.Range("I74").Value = 1665.8182
$ws.Range("J74").Value = 2401.9412
$ws.Range("K74").Value = 1665.8182
$ws.Range("L74").Value = 2401.9412
$ws.Range("M74").Value = -791.8181999999999
$ws.Range("N74").Value = -4149.9412

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2112.75
$ws.Range("I77").Value = 1665.8182
$ws.Range("J77").Value = 2401.9412
$ws.Range("K77").Value = 8329.091
$ws.Range("L77").Value = 12009.706
$ws.Range("M77").Value = -3961.091
$ws.Range("N77").Value = -20745.706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 30238.715
$ws.Range("I102").Value = 1930
$ws.Range("J102").Value = 101010.5
$ws.Range("K102").Value = 1930
$ws.Range("L102").Value = 101010.5
$ws.Range("M102").Value = -308
$ws.Range("N102").Value = -104254.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 12822698
$ws.Range("I132").Value = 16130897
$ws.Range("J132").Value = 3429.5
$ws.Range("K132").Value = 48392691
$ws.Range("L132").Value = 10288.5
$ws.Range("M132").Value = -48390161
$ws.Range("N132").Value = -15348.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 41999.5
$ws.Range("J60").Value = 79999
$ws.Range("L60").Value = 79999
$ws.Range("N60").Value = -81197

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2182.75
$ws.Range("I105").Value = 1988.6364
$ws.Range("K105").Value = 1988.6364
$ws.Range("M105").Value = -241.6364000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 44633.41
$ws.Range("J132").Value = 44633.41
$ws.Range("L132").Value = 44633.41
$ws.Range("N132").Value = -54753.41

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3429.5715
$ws.Range("I134").Value = 2995
$ws.Range("J134").Value = 4516
$ws.Range("K134").Value = 8985
$ws.Range("L134").Value = 13548
$ws.Range("M134").Value = -6450
$ws.Range("N134").Value = -18618

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5245.262
$ws.Range("I31").Value = 2378.3
$ws.Range("J31").Value = 8019.7417
$ws.Range("K31").Value = 2378.3
$ws.Range("L31").Value = 8019.7417
$ws.Range("M31").Value = -2083.3
$ws.Range("N31").Value = -8609.741699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5245.262
$ws.Range("I34").Value = 2378.3
$ws.Range("J34").Value = 8019.7417
$ws.Range("K34").Value = 2378.3
$ws.Range("L34").Value = 8019.7417
$ws.Range("M34").Value = -2176.3
$ws.Range("N34").Value = -8423.741699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2784.6924
$ws.Range("I105").Value = 2669
$ws.Range("J105").Value = 3170.3333
$ws.Range("K105").Value = 2669
$ws.Range("L105").Value = 3170.3333
$ws.Range("M105").Value = -922
$ws.Range("N105").Value = -6664.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 38955.473
$ws.Range("I132").Value = 1538.2858
$ws.Range("J132").Value = 85176.7
$ws.Range("K132").Value = 4614.857400000001
$ws.Range("L132").Value = 255530.1
$ws.Range("M132").Value = -2084.857400000001
$ws.Range("N132").Value = -260590.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 282219.1
$ws.Range("I121").Value = 295.6
$ws.Range("J121").Value = 517155.34
$ws.Range("K121").Value = 886.8000000000001
$ws.Range("L121").Value = 1551466.02
$ws.Range("M121").Value = 423.1999999999999
$ws.Range("N121").Value = -1554086.02

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2288.229
$ws.Range("J131").Value = 1112.7821
$ws.Range("L131").Value = 3338.3463
$ws.Range("N131").Value = -13418.3463

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 66870896
$ws.Range("I134").Value = 77157760
$ws.Range("J134").Value = 6250
$ws.Range("K134").Value = 231473280
$ws.Range("L134").Value = 18750
$ws.Range("M134").Value = -231468210
$ws.Range("N134").Value = -28890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 7466.375
$ws.Range("I97").Value = 893.75
$ws.Range("J97").Value = 14039
$ws.Range("K97").Value = 893.75
$ws.Range("L97").Value = 14039
$ws.Range("M97").Value = -397.75
$ws.Range("N97").Value = -15031

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2339.3489
$ws.Range("I132").Value = 1750
$ws.Range("J132").Value = 3699.3845
$ws.Range("K132").Value = 5250
$ws.Range("L132").Value = 11098.1535
$ws.Range("M132").Value = -2720
$ws.Range("N132").Value = -16158.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1766.5555
$ws.Range("I100").Value = 1612.375
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1612.375
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1071.375
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1598.2727
$ws.Range("I81").Value = 1558.1
$ws.Range("K81").Value = 3116.2
$ws.Range("M81").Value = -2055.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1598.2727
$ws.Range("I84").Value = 1558.1
$ws.Range("K84").Value = 15581
$ws.Range("M84").Value = -10277

Write-Output "edits applied"
